$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2355
$ws.Range("C2").Value = 0.571

$ws.Range("B3").Value = 2041
$ws.Range("C3").Value = 0.583

$ws.Range("C5").Value = 0.667

$ws.Range("C6").Value = 0.308

$ws.Range("C7").Value = 0.261

$ws.Range("B9").Value = 745

$ws.Range("C10").Value = 0.308

$ws.Range("C11").Value = 0.545

$ws.Range("C12").Value = 0.375

$ws.Range("C13").Value = 0.471

$ws.Range("C15").Value = 0.333

$ws.Range("B16").Value = 708
